# "test add new tasks completed"
# Insert a new worksheet "doAddTasks" right after "doSignIn" containing a
# simple task list (title + 3 test tasks), and make it the active sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately before the current 2nd sheet (dnu), i.e.
# right after "doSignIn", so the final order is:
#   doSignIn, doAddTasks, dnu, doAddContact, doFilterSearchInContacts
$doSignIn = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $doSignIn)
$ws.Name = "doAddTasks"

# Populate the task list
$ws.Range("A1").Value = "title"
$ws.Range("A2").Value = "Test add task1"
$ws.Range("A3").Value = "Test add task2"
$ws.Range("A4").Value = "Test add task3"

# Match the column width used for the task title column
$ws.Columns.Item(1).ColumnWidth = 12.5

# Leave the selection on the new sheet where it was left in the saved file
$ws.Range("D12").Select()
